# The underlying edit swaps the full contents of row 4 and row 5 (every
# column A..AY) - the two observation records traded places while staying
# on rows 4 and 5. Most columns (dates, booleans, location names, etc.) hold
# the exact same value in both rows already, so only the columns whose
# content actually differs between row 4 and row 5 need to be written here;
# leaving the rest alone keeps their original cell formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns with a real value in both rows that simply trade places.
$swapCols = "A", "B", "E", "F", "G", "H", "Q", "R", "AM", "AO"

foreach ($col in $swapCols) {
    $c4 = $ws.Range($col + "4")
    $c5 = $ws.Range($col + "5")
    $v4 = $c4.Value2
    $v5 = $c5.Value2
    $c4.Value2 = $v5
    $c5.Value2 = $v4
}

# M / AC hold text only on row 4 (row 5's cell is blank); after the swap the
# text moves to row 5 and row 4 becomes blank.
$m4 = $ws.Range("M4").Value2
$ws.Range("M4").ClearContents()
$ws.Range("M5").Value2 = $m4

$ac4 = $ws.Range("AC4").Value2
$ws.Range("AC4").ClearContents()
$ws.Range("AC5").Value2 = $ac4

# J / L / AF are blank on one row and simply absent on the other; after the
# swap the blank placeholder moves to the opposite row.
$ws.Range("J5").ClearContents()
$ws.Range("J4").Value2 = ""

$ws.Range("L4").ClearContents()
$ws.Range("L5").Value2 = ""

$ws.Range("AF4").ClearContents()
$ws.Range("AF5").Value2 = ""
